$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $c = $ws.Range($range)
    $origStyle = $c.Style
    $c.NumberFormat = "@"
    $c.Value = $value
    $c.Style = $origStyle
}

# Row 2
Set-TextValue "D2" "29.199.68"
$ws.Range("E2").Value = "  -0.74%  "

# Row 3
Set-TextValue "D3" "1.828.15"
$ws.Range("E3").Value = "  -0.79%  "

# Row 4
Set-TextValue "D4" "1.003"
$ws.Range("E4").Value = "  +0.34%  "

# Row 5
Set-TextValue "D5" "234.07"
$ws.Range("E5").Value = "  -2.24%  "

# Row 6
Set-TextValue "D6" "0.5997"
$ws.Range("E6").Value = "  -4.22%  "

# Row 7
$ws.Range("E7").Value = "  +0.24%  "

# Row 8
Set-TextValue "D8" "0.06977"
$ws.Range("E8").Value = "  -5.68%  "

# Row 9
Set-TextValue "D9" "0.2757"
$ws.Range("E9").Value = "  -4.68%  "

# Row 10
Set-TextValue "D10" "23.27"
$ws.Range("E10").Value = "  -6.59%  "

# Row 11
Set-TextValue "D11" "0.07609"
$ws.Range("E11").Value = "  -1.41%  "

# Row 12
Set-TextValue "D12" "1.842.98"
$ws.Range("E12").Value = "  +0.03%  "

# Row 13
Set-TextValue "D13" "4.760"
$ws.Range("E13").Value = "  -4.13%  "

# Row 14
Set-TextValue "D14" "0.6270"
$ws.Range("E14").Value = "  -6.63%  "

# Row 15
Set-TextValue "D15" "0.000009657"
$ws.Range("E15").Value = "  -6.90%  "

# Row 16
Set-TextValue "D16" "78.35"
$ws.Range("E16").Value = "  -4.21%  "

# Row 17
Set-TextValue "D17" "28.786.51"
$ws.Range("E17").Value = "  -2.14%  "

# Row 18
Set-TextValue "D18" "5.709"
$ws.Range("E18").Value = "  -8.89%  "

# Row 19
Set-TextValue "D19" "221.11"
$ws.Range("E19").Value = "  -5.62%  "

# Row 20
$ws.Range("E20").Value = "  +0.29%  "

# Row 21
$ws.Range("E21").Value = "  -6.30%  "

# Row 22
Set-TextValue "D22" "6.857"
$ws.Range("E22").Value = "  -6.16%  "

# Row 23
Set-TextValue "D23" "1.004"
$ws.Range("E23").Value = "  +0.27%  "

# Row 24
Set-TextValue "D24" "155.44"
$ws.Range("E24").Value = "  -0.96%  "

# Row 25
Set-TextValue "D25" "7.968"
$ws.Range("E25").Value = "  -6.00%  "

# Row 26
$ws.Range("E26").Value = "  -4.12%  "

# Row 27
Set-TextValue "D27" "16.54"
$ws.Range("E27").Value = "  -4.50%  "

# Row 28
Set-TextValue "D28" "0.06521"
$ws.Range("E28").Value = "  -10.24%  "

# Row 29
Set-TextValue "D29" "1.452"
$ws.Range("E29").Value = "  -2.30%  "

# Row 30
Set-TextValue "D30" "1.436"
$ws.Range("E30").Value = "  -2.96%  "

# Row 31
Set-TextValue "D31" "3.836"
$ws.Range("E31").Value = "  -4.93%  "

# Row 32
Set-TextValue "D32" "3.762"
$ws.Range("E32").Value = "  -6.70%  "

# Row 33
Set-TextValue "D33" "1.095"
$ws.Range("E33").Value = "  -5.97%  "

# Row 34
Set-TextValue "D34" "1.723"
$ws.Range("E34").Value = "  -5.18%  "

# Row 35
$ws.Range("E35").Value = "  -9.24%  "

# Row 36
Set-TextValue "D36" "2.536"
$ws.Range("E36").Value = "  -1.41%  "

# Row 37
$ws.Range("E37").Value = "  -2.08%  "

# Row 38
Set-TextValue "D38" "0.01748"
$ws.Range("E38").Value = "  -4.93%  "

# Row 39
Set-TextValue "D39" "6.529"
$ws.Range("E39").Value = "  -3.95%  "

# Row 40
Set-TextValue "D40" "1.173.52"
$ws.Range("E40").Value = "  -4.89%  "

# Row 41
Set-TextValue "D41" "0.8967"
$ws.Range("E41").Value = "  -6.22%  "

# Row 42
$ws.Range("E42").Value = "  +0.28%  "

# Row 43
Set-TextValue "D43" "1.982.57"
$ws.Range("E43").Value = "  -0.67%  "

# Row 44
Set-TextValue "D44" "100.50"
$ws.Range("E44").Value = "  -0.71%  "

# Row 45
Set-TextValue "D45" "62.11"
$ws.Range("E45").Value = "  -4.93%  "

# Row 46
Set-TextValue "D46" "0.00000000112"
$ws.Range("E46").Value = "  -3.80%  "

# Row 47
$ws.Range("B47").Value = "Cronos"
$ws.Range("C47").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
Set-TextValue "D47" "0.05597"
$ws.Range("E47").Value = "  -1.17%  "

# Row 48
$ws.Range("B48").Value = "RenderToken"
$ws.Range("C48").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
Set-TextValue "D48" "1.592"
$ws.Range("E48").Value = "  -6.43%  "

# Row 49
$ws.Range("B49").Value = "EnergySwap"
$ws.Range("C49").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-TextValue "D49" "8.470"
$ws.Range("E49").Value = "  -5.35%  "

# Row 50
Set-TextValue "D50" "0.4552"
$ws.Range("E50").Value = "  -0.47%  "

# Row 51
Set-TextValue "D51" "0.3643"
$ws.Range("E51").Value = "  -6.10%  "
